# Generate Report for Handoff
# Updates the localization-status report: sets Priority = "ht" for the
# "Ready for handoff" rows on the zh-cn and de-de sheets, and bumps the
# "Latest Handoff Datetime" timestamp for those same rows.

$wb = $excel.ActiveWorkbook

$rows = @(7, 8, 9, 11, 13, 14)

$zhcn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $zhcn.Range("E$r").Value = "ht"
    $zhcn.Range("H$r").Value = "2016-08-27 04:20:58"
}

$dede = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $dede.Range("E$r").Value = "ht"
    $dede.Range("H$r").Value = "2016-08-27 04:21:07"
}
